$p = $ppt.ActivePresentation

# --- 1) Slide 15: reorder the terminate-sequence animation effects ---
#
# Before:
#   click   -> show shape 27 (style.visibility = visible)
#   w/prev  -> show shape 30
#   click   -> show shape 33   (its own separate build step)
#
# After:
#   click   -> show shape 33   (now fires first, on its own)
#   click   -> show shape 27
#   w/prev  -> show shape 30
#
# (MoveTo/MoveBefore/MoveAfter/Index-assignment are no-ops in this host, so
# the reorder is done by deleting the three affected effects and re-adding
# them via AddEffect in the desired order/trigger combination.)

$s15 = $p.Slides.Item(15)
$seq = $s15.TimeLine.MainSequence

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

$shape27 = Get-ShapeById $s15 27
$shape30 = Get-ShapeById $s15 30
$shape33 = Get-ShapeById $s15 33

# Locate the three effects at the tail of the main sequence (last occurrence
# of each shape id == the final "terminate" build group on this slide).
$idx27 = -1
$idx30 = -1
$idx33 = -1
for ($i = $seq.Count; $i -ge 1; $i--) {
    $eff = $seq.Item($i)
    if ($eff.Shape.Id -eq 27 -and $idx27 -eq -1) { $idx27 = $i }
    if ($eff.Shape.Id -eq 30 -and $idx30 -eq -1) { $idx30 = $i }
    if ($eff.Shape.Id -eq 33 -and $idx33 -eq -1) { $idx33 = $i }
}

# Delete highest index first so earlier indices stay valid.
$indices = @($idx27, $idx30, $idx33) | Sort-Object -Descending
foreach ($ix in $indices) {
    $seq.Item($ix).Delete()
}

# Re-add in the new order: click->33, click->27, withPrevious->30.
$msoAnimEffectAppear = 1
$msoAnimTriggerOnPageClick = 1
$msoAnimTriggerWithPrevious = 2

$seq.AddEffect($shape33, $msoAnimEffectAppear, 0, $msoAnimTriggerOnPageClick) | Out-Null
$seq.AddEffect($shape27, $msoAnimEffectAppear, 0, $msoAnimTriggerOnPageClick) | Out-Null
$seq.AddEffect($shape30, $msoAnimEffectAppear, 0, $msoAnimTriggerWithPrevious) | Out-Null

# --- 2) Slide 3: shorten the bullet text "övning/demo" -> "demo" ---
$s3 = $p.Slides.Item(3)
$shape = $s3.Shapes.Item(2)
$shape.TextFrame.TextRange.Replace("övning/demo", "demo", 0, $false, $false) | Out-Null
